# Implemented sonar selection function. Added pins used to pinlist.
#
# Row 7 gains five new pin/module entries (columns D, G, H, I, J) describing
# the newly-wired sonar analog/digital mux-enable and mux-select GPIOs.
# The row is slightly shorter (13.8pt) and several columns are widened to
# fit the new, longer module descriptions. The active selection is moved
# to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New pin-list entries on row 7 -----------------------------------
$ws.Range("D7").Value = "GPIO (Sonar analog/power mux enable)"
$ws.Range("G7").Value = "GPIO (Sonar digital mux enable)"
$ws.Range("H7").Value = "GPIO (Sonar mux select 0)"
$ws.Range("I7").Value = "GPIO (Sonar mux select 1)"
$ws.Range("J7").Value = "GPIO (Sonar mux select 2)"

# Row 7 height tightened slightly after the edit
$ws.Rows(7).RowHeight = 13.8

# --- Column widths adjusted to fit the new content --------------------
$ws.Columns(7).ColumnWidth = 29.7287449392713
$ws.Columns(9).ColumnWidth = 22.8178137651822
$ws.Columns(10).ColumnWidth = 23.080971659919
$ws.Columns(11).ColumnWidth = 16.165991902834

# --- Move active selection to B4 ---------------------------------------
[void]$ws.Range("B4").Select()
